$wb = $excel.ActiveWorkbook

# --- Summary sheet ---
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1199.25
$summary.Range("B4").Value = -0.75
$summary.Range("B5").Value = -0.88
$summary.Range("B6").Value = 17
$summary.Range("B8").Value = 11
$summary.Range("B9").Value = 29.41

# --- Strategy Status sheet (MarketMaking row, row 4) ---
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 99.25
$status.Range("D4").Value = 17
$status.Range("E4").Value = -0.75
$status.Range("F4").Value = -0.75
$status.Range("G4").Value = 29.41

# --- New trade row (#17) appended to "All Trades" and "MarketMaking" sheets ---
$newRow = @(17, "2026-02-17", "13:17:41", "MarketMaking", "UP", 0.23, 0.2, "CLOSED", -13.0435, -0.03, 99.25, 0, 0, 0.6, "Normal spread capture: 19600 bps", "early_exit", 0.13)

foreach ($sheetName in @("All Trades", "MarketMaking")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $r = 18
    for ($i = 0; $i -lt $newRow.Length; $i++) {
        $c = $i + 1
        $ws.Cells.Item($r, $c).Value = $newRow[$i]
    }
}
